# Update numeric cell values to reflect refreshed algorithm output (Update Name of Algo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E6").Value = 16.37030000000001
$ws.Range("D7").Value = -7.129500000000005
$ws.Range("B9").Value = 6.319399999999997
$ws.Range("D12").Value = -6.885500000000003
$ws.Range("B13").Value = 6.375699999999999
$ws.Range("D14").Value = -7.977400000000005
$ws.Range("E15").Value = 16.23260000000002
$ws.Range("B16").Value = 4.573000000000001
$ws.Range("B18").Value = 7.401299999999997
$ws.Range("D19").Value = -8.050099999999997
$ws.Range("B20").Value = 8.818500000000002
$ws.Range("B26").Value = 4.486300000000008
$ws.Range("D26").Value = -8.963899999999997
$ws.Range("B27").Value = 5.359700000000006
$ws.Range("D27").Value = -9.157399999999996
$ws.Range("E28").Value = 16.49149999999999
$ws.Range("B29").Value = 5.084000000000001
$ws.Range("D29").Value = -7.241499999999995
$ws.Range("E33").Value = 17.10830000000002
$ws.Range("B35").Value = 8.590700000000002
$ws.Range("E35").Value = 16.42820000000001
$ws.Range("B36").Value = 9.081400000000007
$ws.Range("D37").Value = -7.4966
$ws.Range("D38").Value = -7.166300000000001
$ws.Range("E38").Value = 17.2895
$ws.Range("E43").Value = 17.1698
$ws.Range("E44").Value = 16.5463
$ws.Range("B45").Value = 4.845900000000006
$ws.Range("E45").Value = 16.7923
$ws.Range("D47").Value = -7.058600000000002
$ws.Range("E47").Value = 17.21000000000002
$ws.Range("D51").Value = -8.756300000000003
$ws.Range("E51").Value = 16.42769999999999
$ws.Range("D52").Value = -7.2865
$ws.Range("E54").Value = 16.5727
$ws.Range("B55").Value = 6.989199999999996
$ws.Range("D55").Value = -8.116500000000002
$ws.Range("B57").Value = 5.151499999999996
$ws.Range("E57").Value = 16.28749999999999
$ws.Range("E62").Value = 16.1846
$ws.Range("E63").Value = 18.24750000000002
$ws.Range("E67").Value = 17.24800000000002
$ws.Range("B69").Value = 5.200899999999999
$ws.Range("D69").Value = -7.005499999999994
$ws.Range("D70").Value = -7.696100000000005
$ws.Range("E70").Value = 16.99340000000001
$ws.Range("B76").Value = 5.245300000000001
$ws.Range("D76").Value = -8.212499999999995
$ws.Range("B78").Value = 9.928199999999999
$ws.Range("D81").Value = -8.484299999999998
$ws.Range("E81").Value = 16.81609999999998
$ws.Range("B82").Value = 6.706400000000002
$ws.Range("B83").Value = 5.1874
$ws.Range("D83").Value = -9.188999999999989
$ws.Range("E88").Value = 16.2933
$ws.Range("B93").Value = 6.672399999999997
$ws.Range("D94").Value = -7.205899999999999
$ws.Range("E96").Value = 16.24339999999999
$ws.Range("B97").Value = 6.262499999999999
$ws.Range("E99").Value = 16.64550000000001
$ws.Range("D100").Value = -8.577299999999997
$ws.Range("D102").Value = -7.5674
